$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("D2").Value = [double]"1.031401645158516E-17"
$ws.Range("E2").Value = [double]"1.031401645158516E-17"

# Row 3
$ws.Range("D3").Value = [double]"3.10493676625513E-09"
$ws.Range("E3").Value = [double]"3.10493676625513E-09"

# Row 4
$ws.Range("D4").Value = [double]"1.073771617599114E-156"
$ws.Range("E4").Value = [double]"1.073771617599114E-156"

# Row 5
$ws.Range("D5").Value = [double]"1.805185292894419E-71"
$ws.Range("E5").Value = [double]"1.805185292894419E-71"

# Row 6
$ws.Range("D6").Value = [double]"3.230884778053187E-21"
$ws.Range("E6").Value = [double]"3.230884778053187E-21"

# Row 7
$ws.Range("C7").Value = $true
$ws.Range("D7").Value = [double]"0.9999999918666663"
$ws.Range("E7").Value = [double]"8.133333673221443E-09"

# Row 8
$ws.Range("D8").Value = [double]"8.087045002740364E-13"
$ws.Range("E8").Value = [double]"0.9999999999991913"

# Row 10
$ws.Range("D10").Value = [double]"5.924623545835132E-18"

# Row 11
$ws.Range("D11").Value = [double]"0.9999761750033686"
$ws.Range("E11").Value = [double]"2.382499663144966E-05"
$ws.Range("F11").Value = [double]"176.2105407714844"
$ws.Range("G11").Value = 0.7
